# Insert two new data rows at the top of the "Limón" price block (rows 789-790),
# pushing the existing rows 789..877 down to 791..879.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 789 (shifts everything below down by 2)
$ws.Rows("789:790").Insert()

# ---- New row 789 ----
$ws.Range("A789").Value = 7
$ws.Range("B789").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C789").Value = "Ñuble"
$ws.Range("D789").Value = 44918
$ws.Range("E789").Value = 16
$ws.Range("F789").Value = "Fruta"
$ws.Range("G789").Value = 100102
$ws.Range("H789").Value = "Cítricos"
$ws.Range("I789").Value = 100102003
$ws.Range("J789").Value = "Limón"
$ws.Range("K789").Value = "Sin especificar"
$ws.Range("L789").Value = "1a amarillo"
$ws.Range("M789").Value = 160
$ws.Range("N789").Value = 13000
$ws.Range("O789").Value = 14000
$ws.Range("P789").Value = 13500
$ws.Range("Q789").Value = "`$/malla 16 kilos"
$ws.Range("R789").Value = "Región de O'Higgins"
$ws.Range("S789").Value = 844
$ws.Range("T789").Value = 16

# ---- New row 790 ----
$ws.Range("A790").Value = 7
$ws.Range("B790").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C790").Value = "Ñuble"
$ws.Range("D790").Value = 44918
$ws.Range("E790").Value = 16
$ws.Range("F790").Value = "Fruta"
$ws.Range("G790").Value = 100102
$ws.Range("H790").Value = "Cítricos"
$ws.Range("I790").Value = 100102003
$ws.Range("J790").Value = "Limón"
$ws.Range("K790").Value = "Sin especificar"
$ws.Range("L790").Value = "2a amarillo"
$ws.Range("M790").Value = 80
$ws.Range("N790").Value = 12000
$ws.Range("O790").Value = 12000
$ws.Range("P790").Value = 12000
$ws.Range("Q790").Value = "`$/malla 16 kilos"
$ws.Range("R790").Value = "Región de O'Higgins"
$ws.Range("S790").Value = 750
$ws.Range("T790").Value = 16

# Make sure the date cells keep the date number format used throughout column D
$ws.Range("D789:D790").NumberFormat = $ws.Range("D788").NumberFormat
